$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Patient name fields (row 6)
$ws.Range("A6").Value = "FIGUEROA"
$ws.Range("C6").Value = "VILLEGAS"
$ws.Range("E6").Value = "LISBETH"
$ws.Range("G6").Value = "CARLOTA"
$ws.Range("I6").Value = "/201761828"

# Birth date / age / place of birth (row 12)
$ws.Range("A12").Value = "1961-08-15"
$ws.Range("F12").Value = "56"
$ws.Range("H12").Value = "GUATEMLA"

# Occupation / nationality / no. cedula (row 14)
$ws.Range("D14").Value = "AMA DE CASA"
$ws.Range("F14").Value = "GUATEMALTECA"
$ws.Range("H14").Value = "NO PRESENTO"

# Emergency contact (row 20)
$ws.Range("A20").Value = "HERNAN RIOS"
$ws.Range("F20").Value = "ESPOSO"
$ws.Range("H20").Value = "AV. ELENA B 15-11 ZONA 1"
$ws.Range("J20").Value = "42097131"

# Reprint date and time (row 24)
$ws.Range("A24").Value = "24/10/2017"
$ws.Range("C24").Value = "15:25:49"
